$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2").Value = "68.617.66"
$ws.Range("E2").Value = "  +0.01%  "

# Row 3: 'Ethereum'
$ws.Range("D3").Value = "3.909.14"
$ws.Range("E3").Value = "  +0.05%  "

# Row 4: 'TetherUSD'
$ws.Range("E4").Value = "  -0.04%  "

# Row 5: 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "

# Row 6: 'Solana'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.83%  "

# Row 7: 'LidoStakedEther'
$ws.Range("D7").Value = "3.908.95"
$ws.Range("E7").Value = "  +0.07%  "

# Row 8: 'USDC'
$ws.Range("E8").Value = "  +0.00%  "

# Row 9: 'XRP'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.62%  "

# Row 10: 'Dogecoin'
$ws.Range("E10").Value = "  -0.61%  "

# Row 11: 'Toncoin'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "

# Row 12: 'Cardano'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "

# Row 13: 'ShibaInu'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.63%  "

# Row 14: 'Avalanche'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.39%  "

# Row 15: 'WrappedliquidstakedEther2.0'
$ws.Range("D15").Value = "4.559.93"
$ws.Range("E15").Value = "  -0.07%  "

# Row 16: 'WrappedEther'
$ws.Range("D16").Value = "3.908.30"
$ws.Range("E16").Value = "  -0.13%  "

# Row 17: 'WrappedBTC'
$ws.Range("D17").Value = "68.530.71"
$ws.Range("E17").Value = "  -0.27%  "

# Row 18: 'Chainlink'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.13%  "

# Row 19: 'Polkadot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.63%  "

# Row 20: 'TRON'
$ws.Range("E20").Value = "  +0.39%  "

# Row 21: 'Uniswap'
$ws.Range("E21").Value = "  -2.02%  "

# Row 22: 'BitcoinCash'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "474.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.84%  "

# Row 23: 'Polygon'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.742"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.06%  "

# Row 24: 'PEPE'
$ws.Range("E24").Value = "  +0.23%  "

# Row 25: 'Litecoin'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.97%  "

# Row 26: 'Fetch.AI'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.33%  "

# Row 27: 'InternetComputer(DFINITY)'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.44%  "

# Row 28: 'Dai' -> 'RenderToken'
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.11%  "

# Row 29: 'RenderToken' -> 'Dai'
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "

# Row 30: 'PancakeSwap'
$ws.Range("E30").Value = "  +1.45%  "

# Row 31: 'WrappedeETH'
$ws.Range("D31").Value = "4.053.89"
$ws.Range("E31").Value = "  -0.18%  "

# Row 32: 'NEARProtocol'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.31%  "

# Row 33: 'EthereumClassic'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.91%  "

# Row 34: 'ImmutableX'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.18%  "

# Row 35: 'Aptos'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.62%  "

# Row 36: 'RenzoRestakedETH'
$ws.Range("D36").Value = "3.876.25"
$ws.Range("E36").Value = "  +0.30%  "

# Row 37: 'Hedera'
$ws.Range("E37").Value = "  -1.29%  "

# Row 38: 'dogwifhat'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.20%  "

# Row 39: 'Kaspa' -> 'Mantle'
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.18%  "

# Row 40: 'Mantle' -> 'Kaspa'
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.142"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.00%  "

# Row 41: 'Filecoin'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.15%  "

# Row 42: 'FirstDigitalUSD'
$ws.Range("E42").Value = "  -0.06%  "

# Row 43: 'TheGraph'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.315"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.61%  "

# Row 44: 'FLOKI' -> 'Bittensor'
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "429.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.36%  "

# Row 45: 'Stacks'
$ws.Range("E45").Value = "  +0.56%  "

# Row 46: 'Bittensor' -> 'FLOKI'
$ws.Range("B46").Value = "FLOKI"
$ws.Range("C46").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000301"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +13.37%  "

# Row 47: 'USDe' -> 'Cosmos'
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.71%  "

# Row 48: 'Cosmos' -> 'USDe'
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "

# Row 49: 'OKB'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "47.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.99%  "

# Row 50: 'EnergySwap'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.69%  "

# Row 51: 'Monero'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.07%  "

